$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("by Coach")

# Flip the "Started" (column C) Yes/No flags for the rows that changed.
$changes = @{
    2  = "Yes"
    3  = "No"
    4  = "Yes"
    7  = "No"
    8  = "Yes"
    9  = "No"
    14 = "Yes"
    15 = "Yes"
    21 = "No"
    24 = "No"
    29 = "Yes"
    34 = "No"
    42 = "Yes"
    45 = "No"
    51 = "Yes"
    53 = "Yes"
    54 = "No"
    56 = "No"
    58 = "No"
    60 = "Yes"
    64 = "No"
    69 = "Yes"
    82 = "Yes"
    83 = "No"
}

foreach ($row in $changes.Keys) {
    $ws.Range("C$row").Value = $changes[$row]
}

# Scroll the frozen pane down near the bottom of the list, then land the
# selection on C83, matching where the editor ended up.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$ws.Range("C83").Select()
